$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Table 1, row 5 ("16:00 - 17:00" / "Ontwerp document af") - fill in the
#    Start/Eind times and the remark cell, matching how row 4 was filled in.
# ---------------------------------------------------------------------------

function Get-Row5TimeRange {
    $t = $d.Tables.Item(1)
    $cell = $t.Cell(5, 3)
    return $d.Range($cell.Range.Start, $cell.Range.End)
}

# Start: <00:00> -> Start: <20:30>
$r = Get-Row5TimeRange
$r.Find.Execute("00", $false, $false, $false, $false, $false, $true, 1, $false, "20", 1) | Out-Null

$r = Get-Row5TimeRange
$r.Find.Execute("00>", $false, $false, $false, $false, $false, $true, 1, $false, "30>", 1) | Out-Null

# Eind: <00:00> -> Eind: <23:05>
$r = Get-Row5TimeRange
$r.Find.Execute("00", $false, $false, $false, $false, $false, $true, 1, $false, "23", 1) | Out-Null

$r = Get-Row5TimeRange
$r.Find.Execute("00>", $false, $false, $false, $false, $false, $true, 1, $false, "05>", 1) | Out-Null

# Remark cell: was empty, now gets the "2,5 uur ipv 1 uur." remark plus the
# "_GoBack" bookmark that used to sit after "Testmatrix hoefde niet." in the
# row above.
$t = $d.Tables.Item(1)
$remarkCell = $t.Cell(5, 4)
$remarkRange = $remarkCell.Range
$remarkRange.Collapse(0)
$remarkRange.InsertBefore("2,5 uur ipv 1 uur.")
$remarkRange.InsertAfter([char]11)
$remarkRange.InsertAfter("Geen rekening met gui en diagrammen gehouden. ")

$t = $d.Tables.Item(1)
$remarkCell = $t.Cell(5, 4)
$bmRange = $remarkCell.Range
$bmRange.Collapse(0)
$bmRange.MoveEnd(1, -1) | Out-Null
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------------
# 2) "(met id)" - no textual change (markup-only proofing annotation in the
#    source diff), left as-is.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 3) "Use case diagram" / "Use cases" - no textual change either.
# ---------------------------------------------------------------------------
